$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the "Potential Form" (column G) for rows 68-84, which previously
#    held the placeholder text "na". Most of these cells also switch their
#    font from "微軟正黑體" (style s="4") to "Yu Gothic" (style s="3") to
#    match the Japanese text now being displayed; two rows (74 and 77) keep
#    their original font.
# ---------------------------------------------------------------------------
$potentialForms = @{
    68 = "渡れる"
    69 = "払える"
    70 = "謝れる"
    71 = "困れる"
    72 = "始まれる"
    73 = "終われる"
    74 = "掛れる"
    75 = "押せる"
    76 = "渡せる"
    77 = "返せる"
    78 = "焼ける"
    79 = "とおれる"
    80 = "かよえる"
    81 = "送れる"
    82 = "上がれる"
    83 = "下ろせる"
    84 = "下がれる"
}

# Rows that keep the original (微軟正黑體) font instead of switching to Yu Gothic
$keepOriginalFont = @(74, 77)

foreach ($row in 68..84) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.Value = $potentialForms[$row]
    if ($keepOriginalFont -notcontains $row) {
        $cell.Font.Name = "Yu Gothic"
    }
}

# ---------------------------------------------------------------------------
# 2) Add a brand new verb entry (誘う / "to invite") as row 107, using the
#    same "Yu Gothic" styled font as the rest of the recently-added rows.
# ---------------------------------------------------------------------------
$newRow = @("誘う", "誘って", "誘った", "誘わない", "誘います", "誘おう", "誘える")

for ($col = 1; $col -le 7; $col++) {
    $cell = $ws.Cells.Item(107, $col)
    $cell.Value = $newRow[$col - 1]
    $cell.Font.Name = "Yu Gothic"
}

$ws.Rows.Item(107).RowHeight = 18.75

# ---------------------------------------------------------------------------
# 3) Widen column H (newly touched, though left without data) to roughly
#    match the authored width of 44.28515625 characters.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 43.5

# ---------------------------------------------------------------------------
# 4) Restore the cursor/selection position used when the file was last saved.
# ---------------------------------------------------------------------------
$ws.Range("F106").Select()
